$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# Copy the date cell formatting (style s="8") from the row above, then set the values.
$ws.Range("C22").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B23").Value = "D29"
$ws.Range("C23").Value = 43863
$ws.Range("D23").Value = "Completed lesson 7. Worked on Intel Image classification."

$ws.Range("D23").Select()
